# Änderungen nach DA Übergabe III
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Tabelle1 (sheet1): update D2, clear row 3 values, move selection
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws1.Range("D2").Value = 3000
$ws1.Range("A3:E3").ClearContents()
$ws1.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------
# Tabelle2 (sheet2): update row 2 values, clear row 3, add selection
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Tabelle2")
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = 0
$ws2.Range("D2").Value = 0.0002
$ws2.Range("E2").Value = 0.00005
$ws2.Range("F2").Value = 0.00007
$ws2.Range("G2").Value = 0.009
$ws2.Range("I2").Value = 0.003
$ws2.Range("A3:I3").ClearContents()
$ws2.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------
# Tabelle3 (sheet3): update row 2 values, clear row 3, move selection
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Tabelle3")
$ws3.Range("A2").Value = 30.5
$ws3.Range("B2").Value = 7035
$ws3.Range("C2").Value = 492
$ws3.Range("F2").Value = 1793
$ws3.Range("A3:H3").ClearContents()
$ws3.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------
# Tabelle4 (sheet4): clear row 3 values, move selection to A3
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Tabelle4")
$ws4.Range("A3:B3").ClearContents()
$ws4.Range("A3").Select() | Out-Null

# ---------------------------------------------------------------
# Tabelle5 (sheet5): update row 2 values, clear row 3, move selection
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Tabelle5")
$ws5.Range("B2").Value = 20
$ws5.Range("C2").Value = 100
$ws5.Range("A3:B3").ClearContents()
$ws5.Range("C3").Clear()
$ws5.Range("A2").Select() | Out-Null

# Re-select Tabelle2 as the active sheet (tabSelected="1")
$ws2.Activate() | Out-Null
